# ---------------------------------------------------------------
# Week 6 evaluation rubric: populate the newly-added evaluation
# criteria (rows 2-6) with full grading-scale descriptions across
# columns C:F, and a new "Missing" column G. Replaces the old
# placeholder criteria text in column B.
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s0 = @"
Criteria
"@
$s1 = @"
Exceptional
"@
$s2 = @"
Proficient
"@
$s3 = @"
Basic
"@
$s4 = @"
Limited
"@
$s5 = @"
PowerPoint Presentation and Delivery
"@
$s6 = @"
Include
"@
$s7 = @"
y
"@
$s8 = @"
Quality and Relevance of Visualizations
"@
$s9 = @"
Structure and Detail of Data Exploration
"@
$s10 = @"
Y
"@
$s11 = @"
Missing
"@
$s12 = @"
The PowerPoint presentation is exceptional in design and delivery. Slides are cohesive with a consistent design, using appropriate fonts, colors, and objects. The presentation flows smoothly, is engaging, and keeps the audience's attention. The presenter delivers the content confidently and clearly within the allotted time slot. 
"@
$s13 = @"
The PowerPoint presentation is well- designed and delivered. Slides exhibit a good degree of cohesion in terms of design elements. The presentation is clear and mostly fluent, and the presenter manages to stay within the allotted time. 
"@
$s14 = @"
The PowerPoint presentation lacks some cohesiveness in design elements, making it less visually appealing. The delivery is somewhat choppy but generally understandable, and the presentation mostly fits the allotted time. 
"@
$s15 = @"
The PowerPoint presentation is poorly designed, lacking consistency in fonts, colors, and objects. The delivery is disjointed or difficult to follow, and the presentation significantly exceeds or falls short of the allotted time. 
"@
$s16 = @"
Missing or no work was submitted.
"@
$s17 = @"
The data exploration steps are not clearly outlined or missing. The process to explore different use cases for the business problem at hand is not clearly explained.
"@
$s18 = @"
Tension Development and Persuasiveness of  Argumentation in Data Story
"@
$s19 = @"
Clarity of Big Idea / Primary Use Case for a Non-Technical Audience
"@
$s20 = @"
Visually through smart use of text placement, spacing, style, visuals and audibly through speaker's / speakers' tone, pace, and other oratory techniques the audience is guided effectively and convincingly to the final selected use case.
"@
$s21 = @"
The final selected use case is clearly separated from other possible use cases and outlined at the end of the presentation. Verbal and visual presentation are aligned well to communicate the Big Idea / primary use case effectively.
"@
$s22 = @"
The final selected use case is mentioned but could have been better separated from other use cases mentioned.
"@
$s23 = @"
It is not clear from the presentation what the final selected use case is. 
"@
$s24 = @"
The visualizations stand out in terms of quality, style, clarity. They are able to help focus the key messages in the data story. Color, graph choice, labeling, descriptions are thoughtfully and effectively used to create high-quality professional visualizations. With the three visualizations the audience is taken on a journey that leads through the essential aspects of the data to persuasively motivate the solution proposed for the business problem at hand.
"@
$s25 = @"
The visualizations used were relevant and informative and logically created a path for the audience to understand the particular complexities of the data and how handling those in logical sequence will lead to an effective solution of the business problem. This is achieved with the help of three well-crafted visuals tied into a data story that concludes with a data science solution proposal of the business problem. Space, color, the choice of graphing style and other elements add effectively to the telling of the data story.
"@
$s26 = @"
Three visualizations were added that were useful to communicate essemtial parts of the data story but visuals or their descriptions in the presentation could have been improved to tell the data story more clearly.
"@
$s27 = @"
Visualizations were included in the slide deck but only with limited relevance to telling the data story of how the data expolration can help solve the bbusiness problem prsented. Key elements in the graph like axis labels are missing or hard to read or generally visualizations asre missing and/or of lower quality.
"@
$s28 = @"
The data exploration illustrates a high level of thought and insight with a feasibility study for use cases clearly outlined and motivated. The logic succession of steps in the data processing is easily understandable.
"@
$s29 = @"
The data exploration is clearly described and well-formulated as a data story. The logic in the data exploration can be understood and conclusions make sense and carry the data story forward to motivate the proposed solution for business problem.
"@
$s30 = @"
Visual and verbal presentation align perfectly to elegantly and convincingly walk the audience through the data story. The reasoning is persuasive and well designed to also convince non-technical business stakeholders.
"@
$s31 = @"
Visual and verbal presentation components integrate well to make the case for the selection of steps in the data exploration process. The motivation is easy to follow and is convincing based on the data and the analysis presented.
"@
$s32 = @"
It is not clear how some or any of the steps in the data exploration process are relevant for the data science use case and the business problem. 
"@
$s33 = @"
The presentation shows all the relevant steps in the data exploration process. But the motivation for the selection and prioritization of steps In the data exploration process is not clear.
"@
$s34 = @"
A data exploration process is described with some explanation of how certain aspects in the data will affect the use case discoery process. But the description of feature importance and feature correlations and correlation to the target feature is not clear and it is hard to follow a data story that concludes with a solution towards the business problem.
"@

# --- header row ---
$ws.Range("A1").Value = $s6
$ws.Range("B1").Value = $s0
$ws.Range("C1").Value = $s1
$ws.Range("D1").Value = $s2
$ws.Range("E1").Value = $s3
$ws.Range("F1").Value = $s4
$ws.Range("G1").Value = $s11

# --- row 2: Clarity of Big Idea / Primary Use Case for a Non-Technical Audience ---
$ws.Range("A2").Value = $s7
$ws.Range("B2").Value = $s19
$ws.Range("C2").Value = $s20
$ws.Range("D2").Value = $s21
$ws.Range("E2").Value = $s22
$ws.Range("F2").Value = $s23
$ws.Range("G2").Value = $s16

# --- row 3: Tension Development and Persuasiveness of Argumentation in Data Story ---
$ws.Range("A3").Value = $s10
$ws.Range("B3").Value = $s18
$ws.Range("C3").Value = $s30
$ws.Range("D3").Value = $s31
$ws.Range("E3").Value = $s33
$ws.Range("F3").Value = $s32
$ws.Range("G3").Value = $s16

# --- row 4: Structure and Detail of Data Exploration ---
$ws.Range("A4").Value = $s7
$ws.Range("B4").Value = $s9
$ws.Range("C4").Value = $s28
$ws.Range("D4").Value = $s29
$ws.Range("E4").Value = $s34
$ws.Range("F4").Value = $s17
$ws.Range("G4").Value = $s16

# --- row 5: Quality and Relevance of Visualizations ---
$ws.Range("A5").Value = $s7
$ws.Range("B5").Value = $s8
$ws.Range("C5").Value = $s24
$ws.Range("D5").Value = $s25
$ws.Range("E5").Value = $s26
$ws.Range("F5").Value = $s27
$ws.Range("G5").Value = $s16

# --- row 6: PowerPoint Presentation and Delivery ---
$ws.Range("A6").Value = $s7
$ws.Range("B6").Value = $s5
$ws.Range("C6").Value = $s12
$ws.Range("D6").Value = $s13
$ws.Range("E6").Value = $s14
$ws.Range("F6").Value = $s15
$ws.Range("G6").Value = $s16

# --- row heights (auto-grown to fit the new wrapped descriptions) ---
$ws.Rows.Item(2).RowHeight = 85
$ws.Rows.Item(3).RowHeight = 85
$ws.Rows.Item(4).RowHeight = 136
$ws.Rows.Item(5).RowHeight = 187
$ws.Rows.Item(6).RowHeight = 119

# --- columns C:F in rows 4 and 6 use a plain (non-theme) Calibri font ---
$ws.Range("C4:F4").Font.Name = "Calibri"
$ws.Range("C6:F6").Font.Name = "Calibri"

# --- restore the author's final selection ---
$ws.Range("F5").Select()
